# Commit: feat(translation): update translation for all languages
#
# 1. Standardise the "English Copy" / "Target Language" / " Bengali" header
#    labels on both data sheets to "English copy" / "Bengali".
# 2. Insert a new blank row at the top of "Delta 2" (its content used to
#    start at row 1; it now starts at row 2, matching "Delta").
# 3. Remove the unused, empty "Sheet10" tab.
# 4. Make "Delta 2" the active sheet / tab, with G3 selected on both sheets.

$wb = $excel.ActiveWorkbook

$deltaSheet  = $wb.Worksheets.Item("Delta")
$delta2Sheet = $wb.Worksheets.Item("Delta 2")

# --- Delta 2: shift everything down by one row, insert blank row 1 ---------
$delta2Sheet.Rows.Item(1).Insert() | Out-Null

# --- Unify the header row text on both sheets ------------------------------
$deltaSheet.Range("E2").Value = "English copy"
$deltaSheet.Range("G2").Value = "Bengali"

$delta2Sheet.Range("E2").Value = "English copy"
$delta2Sheet.Range("G2").Value = "Bengali"

# --- Remove the empty "Sheet10" tab -----------------------------------------
$wb.Worksheets.Item("Sheet10").Delete() | Out-Null

# --- Selection / active tab -------------------------------------------------
$deltaSheet.Range("G3").Select() | Out-Null
$delta2Sheet.Activate() | Out-Null
$delta2Sheet.Range("G3").Select() | Out-Null
